$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a literal text value into a cell without letting Excel
# auto-convert number-looking strings (e.g. "1.001", "262.90") into real
# numbers. Setting NumberFormat "@" first forces text entry; resetting the
# Style to "Normal" afterwards drops the temporary text-format style so the
# cell ends up with the same (default) style index it had before the edit.
function Set-TextCell($addr, $text) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

# --- Rows 2-20: Price (D) / Volume(1h) (E) refresh; Coin (B) / Link (C) unchanged ---
Set-TextCell "D2" "26.631.66"
Set-TextCell "E2" "  +0.11%  "
Set-TextCell "D3" "1.848.68"
Set-TextCell "E3" "  -0.30%  "
Set-TextCell "E4" "  +0.10%  "
Set-TextCell "D5" "262.90"
Set-TextCell "E5" "  -0.81%  "
Set-TextCell "E6" "  +0.10%  "
Set-TextCell "D7" "0.5342"
Set-TextCell "E7" "  +2.11%  "
Set-TextCell "D8" "0.3166"
Set-TextCell "E8" "  -3.78%  "
Set-TextCell "D9" "0.06959"
Set-TextCell "E9" "  +2.11%  "
Set-TextCell "D10" "18.91"
Set-TextCell "E10" "  +0.10%  "
Set-TextCell "D11" "0.7712"
Set-TextCell "E11" "  -0.84%  "
Set-TextCell "D12" "0.07837"
Set-TextCell "E12" "  +0.88%  "
Set-TextCell "D13" "1.858.93"
Set-TextCell "E13" "  +0.30%  "
Set-TextCell "D14" "89.65"
Set-TextCell "E14" "  +1.10%  "
Set-TextCell "D15" "5.047"
Set-TextCell "E15" "  +0.32%  "
Set-TextCell "D16" "14.13"
Set-TextCell "E16" "  +0.81%  "
Set-TextCell "D17" "1.001"
Set-TextCell "E17" "  +0.13%  "
Set-TextCell "D18" "0.000007977"
Set-TextCell "E18" "  -0.01%  "
Set-TextCell "E19" "  +0.20%  "
Set-TextCell "D20" "26.646.27"
Set-TextCell "E20" "  +0.12%  "

# --- Rows 21-51: coin list shifted up one slot (oldest entry dropped, new coin appended at the bottom) ---
Set-TextCell "B21" "Uniswap"
Set-TextCell "C21" "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
Set-TextCell "D21" "4.644"
Set-TextCell "E21" "  -0.03%  "
Set-TextCell "B22" "Chainlink"
Set-TextCell "C22" "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextCell "D22" "6.021"
Set-TextCell "E22" "  +0.39%  "
Set-TextCell "B23" "Cosmos"
Set-TextCell "C23" "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextCell "D23" "9.371"
Set-TextCell "E23" "  -2.14%  "
Set-TextCell "B24" "LidoDAOToken"
Set-TextCell "C24" "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
Set-TextCell "D24" "2.215"
Set-TextCell "E24" "  +0.27%  "
Set-TextCell "B25" "Monero"
Set-TextCell "C25" "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextCell "D25" "141.92"
Set-TextCell "E25" "  -1.87%  "
Set-TextCell "B26" "Toncoin"
Set-TextCell "C26" "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextCell "D26" "1.691"
Set-TextCell "E26" "  +1.57%  "
Set-TextCell "B27" "EthereumClassic"
Set-TextCell "C27" "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextCell "D27" "17.13"
Set-TextCell "E27" "  +0.54%  "
Set-TextCell "B28" "BitcoinCash"
Set-TextCell "C28" "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
Set-TextCell "D28" "111.62"
Set-TextCell "E28" "  -0.48%  "
Set-TextCell "B29" "InternetComputer(DFINITY)"
Set-TextCell "C29" "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextCell "D29" "4.312"
Set-TextCell "E29" "  +2.39%  "
Set-TextCell "B30" "Stellar"
Set-TextCell "C30" "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextCell "D30" "0.08782"
Set-TextCell "E30" "  +0.14%  "
Set-TextCell "B31" "Filecoin"
Set-TextCell "C31" "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextCell "D31" "4.109"
Set-TextCell "E31" "  -1.23%  "
Set-TextCell "B32" "Hedera"
Set-TextCell "C32" "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextCell "D32" "0.04856"
Set-TextCell "E32" "  +0.18%  "
Set-TextCell "B33" "ImmutableX"
Set-TextCell "C33" "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextCell "D33" "0.7373"
Set-TextCell "E33" "  +2.74%  "
Set-TextCell "B34" "ARBITRUM"
Set-TextCell "C34" "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextCell "D34" "1.140"
Set-TextCell "E34" "  -0.08%  "
Set-TextCell "B35" "HuobiToken"
Set-TextCell "C35" "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
Set-TextCell "D35" "2.885"
Set-TextCell "E35" "  +1.07%  "
Set-TextCell "B36" "MXToken"
Set-TextCell "C36" "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextCell "D36" "3.106"
Set-TextCell "E36" "  -0.18%  "
Set-TextCell "B37" "RenderToken"
Set-TextCell "C37" "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextCell "D37" "2.361"
Set-TextCell "E37" "  +6.22%  "
Set-TextCell "B38" "VeChain"
Set-TextCell "C38" "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextCell "D38" "0.01738"
Set-TextCell "E38" "  -2.55%  "
Set-TextCell "B39" "TheSandbox"
Set-TextCell "C39" "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
Set-TextCell "D39" "0.4836"
Set-TextCell "E39" "  -1.42%  "
Set-TextCell "B40" "TrustWalletToken"
Set-TextCell "C40" "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextCell "D40" "0.9090"
Set-TextCell "E40" "  +0.22%  "
Set-TextCell "B41" "Quant"
Set-TextCell "C41" "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
Set-TextCell "D41" "109.02"
Set-TextCell "E41" "  -3.27%  "
Set-TextCell "B42" "FraxShare"
Set-TextCell "C42" "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextCell "D42" "5.915"
Set-TextCell "E42" "  -2.93%  "
Set-TextCell "B43" "PaxDollar"
Set-TextCell "C43" "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
Set-TextCell "D43" "1.001"
Set-TextCell "E43" "  +0.13%  "
Set-TextCell "B44" "Aptos"
Set-TextCell "C44" "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextCell "D44" "7.709"
Set-TextCell "E44" "  -0.50%  "
Set-TextCell "B45" "Decentraland"
Set-TextCell "C45" "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
Set-TextCell "D45" "0.4200"
Set-TextCell "E45" "  -0.04%  "
Set-TextCell "B46" "EnergySwap"
Set-TextCell "C46" "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextCell "D46" "9.124"
Set-TextCell "E46" "  +0.00%  "
Set-TextCell "B47" "Algorand"
Set-TextCell "C47" "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextCell "D47" "0.1250"
Set-TextCell "E47" "  +0.55%  "
Set-TextCell "B48" "Elrond"
Set-TextCell "C48" "https://coinranking.com/coin/omwkOTglq+elrond-egld"
Set-TextCell "D48" "35.07"
Set-TextCell "E48" "  -0.32%  "
Set-TextCell "B49" "Cronos"
Set-TextCell "C49" "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextCell "D49" "0.05816"
Set-TextCell "E49" "  -2.09%  "
Set-TextCell "B50" "EOS"
Set-TextCell "C50" "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
Set-TextCell "D50" "0.8976"
Set-TextCell "E50" "  +0.85%  "
Set-TextCell "B51" "Aave"
Set-TextCell "C51" "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextCell "D51" "60.36"
Set-TextCell "E51" "  +0.43%  "
